$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.023.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.249.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.26%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.79%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.523"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.485"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0810"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.63%  "

$ws.Range("E12").Value = "  +0.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.597.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.16%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.246.56"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.778"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.883.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0900"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "

$ws.Range("E21").Value = "  -0.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.94"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.20%  "

$ws.Range("E29").Value = "  +0.76%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.04%  "

$ws.Range("E32").Value = "  +0.09%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.74%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0721"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.75%  "

$ws.Range("E37").Value = "  -0.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.114"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.38%  "

$ws.Range("E39").Value = "  -3.07%  "

$ws.Range("E40").Value = "  -3.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.937.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.08%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0281"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.71%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -9.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.471.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "91.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.88%  "
